# Insert a new "interval_fine" column between the existing "interval" (AF)
# and "area" (AG, which becomes AH) columns, and populate it with
# inside/outside flags for the finer 99.75-100.5 interval check.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column AG (33rd column) - shifts old AG (area) to AH.
$ws.Columns.Item(33).Insert()

# Header for the new column.
$ws.Range("AG1").Value() = "interval_fine"

# Default every data row (2-180) to "inside" ...
$ws.Range("AG2:AG180").Value() = "inside"

# ... then flip the rows that fall outside the finer interval to "outside".
$outsideRows = @(8, 11, 13, 17, 19, 21, 22, 28, 30, 34, 41, 43, 46, 53, 54, 58, 60, 64, 65, 66, 80, 83, 85, 92, 96, 104, 106, 107, 109, 113, 115, 125, 126, 129, 130, 139, 143, 144, 146, 152, 153, 160, 164, 166, 168, 169, 171, 177, 179)

foreach ($r in $outsideRows) {
    $ws.Cells.Item($r, 33).Value() = "outside"
}
